# issue #5: stock data output to json file
#
# The stock ("股票") sheet gains a new "property_category" column (value
# "stock" on every data row), inserted right after "total" and before
# "date" — pushing date/legislator_name/legislator_id one column to the
# right. The quantity text on the first stock row also loses its
# thousands-separator comma ("1’160,000" -> "1’160000").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = $ws.UsedRange.Rows.Count

# Insert a new blank column before the current "date" column (H); this
# shifts date -> I, legislator_name -> J, legislator_id -> K.
$ws.Columns.Item(8).Insert()

# New column header.
$ws.Cells.Item(1, 8).Value = "property_category"

# New column data: every stock row is categorized as "stock".
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}

# Normalize the quantity text for the first stock row (drop the comma).
$ws.Cells.Item(2, 4).Value = "1’160000"
